$wb = $excel.ActiveWorkbook

# --- "category" sheet: add new "hot" (int) and "where" (string) fields ---
$wsCategory = $wb.Worksheets.Item("category")

$wsCategory.Range("F1").Value = "hot"
$wsCategory.Range("G1").Value = "where"

$wsCategory.Range("F2").Value = "热度"
$wsCategory.Range("G2").Value = "位置"

$wsCategory.Range("F3").Value = "int"
$wsCategory.Range("G3").Value = "string"

$wsCategory.Range("F5").Value = 10
$wsCategory.Range("G5").Value = "二楼"

$wsCategory.Range("F8").Value = 1
$wsCategory.Range("G8").Value = "一楼"

# --- "pages" sheet: add new "extra_1" (int) and "extra_2" (string) fields ---
$wsPages = $wb.Worksheets.Item("pages")

$wsPages.Range("D1").Value = "extra_1"
$wsPages.Range("E1").Value = "extra_2"

$wsPages.Range("D2").Value = "额外 1"

$wsPages.Range("D3").Value = "int"
$wsPages.Range("E3").Value = "string"

$wsPages.Range("D5").Value = 0

$wsPages.Range("E7").Value = "我说哪。"

$wsPages.Range("D8").Value = 3

# --- selection / active-tab bookkeeping ---
# "pages" keeps the focus on its last-used cell
$wsPages.Activate() | Out-Null
$wsPages.Range("E8").Select() | Out-Null

# "book" was the previously active tab; it no longer is
$wsBook = $wb.Worksheets.Item("book")
$wsBook.Activate() | Out-Null
$wsBook.Range("H9").Select() | Out-Null

# "category" becomes the active tab, selection parked past the new columns
$wsCategory.Activate() | Out-Null
$wsCategory.Range("G9").Select() | Out-Null
